# OpenData Slovakia Covid Deaths Cumulative - append new daily rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the existing row 44 (2020-11-26), shifting rows
# 44:66 down to 45:67, then populate it with the 2020-11-25 data point.
$ws.Rows.Item(44).Insert()

$ws.Range("A44").Value = 44160
$ws.Range("A44").Style = $ws.Range("A45").Style
$ws.Range("B44").Value = 749
$ws.Range("C44").Value = 196
$ws.Range("D44").Value = 945

# Append the new daily rows for 2020-12-19 through 2020-12-28 after the
# (now shifted) last existing row, 67.
$newRows = @(
  @(44184, 1555, 394, 1949),
  @(44185, 1618, 397, 2015),
  @(44186, 1655, 412, 2067),
  @(44187, 1686, 420, 2106),
  @(44188, 1732, 429, 2161),
  @(44191, 1773, 458, 2231),
  @(44192, 1879, 498, 2377),
  @(44193, 1983, 501, 2484)
)

$r = 68
foreach ($row in $newRows) {
  $ws.Range("A$r").Value = $row[0]
  $ws.Range("A$r").Style = $ws.Range("A67").Style
  $ws.Range("B$r").Value = $row[1]
  $ws.Range("C$r").Value = $row[2]
  $ws.Range("D$r").Value = $row[3]
  $r = $r + 1
}

# Restore the view: scroll position and active cell selection match the
# saved workbook state after the edit.
$ws.Application.ActiveWindow.ScrollRow = 29
$ws.Range("A75").Select()
